$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = "en los logacciones deberia aparecer que terminó el cuarto y al momento de hacer la falta debe informar, tipo de falta, si da lanzamientos y cantidad de faltas de ese jugador"
$ws.Range("A24").Value = "tambien deben aparecer las correcciones"

$ws.Range("A25").Select()
$excel.ActiveWindow.ScrollRow = 7

